$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new record above the current row 168, shifting the
# existing rows 168-230 down to 169-231 (dimension grows to A1:R231).
$ws.Rows.Item(168).EntireRow.Insert()

# Populate the newly inserted row 168 with the new "Papa" price record.
$ws.Range("A168").Value = 4
$ws.Range("B168").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C168").Value = "Los Lagos"
$ws.Range("D168").Value = 44468
$ws.Range("E168").Value = 10
$ws.Range("F168").Value = 100114001
$ws.Range("G168").Value = "Papa"
$ws.Range("H168").Value = "Asterix"
$ws.Range("I168").Value = "1a (guarda)"
$ws.Range("J168").Value = 150
$ws.Range("K168").Value = 8000
$ws.Range("L168").Value = 8000
$ws.Range("M168").Value = 8000
$ws.Range("N168").Value = '$/saco 25 kilos'
$ws.Range("O168").Value = "Provincia de Llanquihue"
$ws.Range("P168").Value = 320
$ws.Range("Q168").Value = 25
$ws.Range("R168").Value = "Hortaliza"
